$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card14")

# Row 17 used to be the last data row in the table; its previously-blank
# B..K cells get filled in with the literal "nan" placeholder text (as
# produced by the pandas/openpyxl export this workbook is generated from)
# now that a new row has been appended after it.
$ws.Range("B17:K17").Value = "nan"

# Row 18 is the newly added service-card event.
# Column A ("card") must stay text "14" like every other row in this
# column, not get auto-converted to a number by Excel.
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "14"
$ws.Range("A18").ClearFormats()

# Columns B..K are intentionally left untouched/blank for this new row.

$ws.Range("L18").Value = "12\8\2025"
$ws.Range("M18").Value = "777 t"
$ws.Range("N18").Value = "تم تغيير زيت الجيربوكس وتغيير جريد 1"
$ws.Range("O18").Value = "تيم العمل"
